# Scheduled-runner refresh of cached market-board figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across all job sheets. Values only - no formulas,
# formatting, or structural changes. A few rows gain/lose trailing M/N cells
# where the cached profit figure became present/absent upstream.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1206.375
$ws.Range("I15").Value = 1206.375
$ws.Range("K15").Value = 3619.125
$ws.Range("M15").Value = -3450.125
$ws.Range("H18").Value = 387.5
$ws.Range("I18").Value = 387.5
$ws.Range("K18").Value = 387.5
$ws.Range("M18").Value = -103.5
$ws.Range("H116").Value = 23250.5
$ws.Range("I116").Value = 7832.6665
$ws.Range("K116").Value = 7832.6665
$ws.Range("M116").Value = -4390.6665
$ws.Range("H125").Value = 15876619.0
$ws.Range("I125").Value = 2000.0
$ws.Range("J125").Value = 18522388.0
$ws.Range("K125").Value = 18000.0
$ws.Range("L125").Value = 166701492.0
$ws.Range("M125").Value = -15540.0
$ws.Range("N125").Value = -166706412.0
$ws.Range("H138").Value = 5767.976
$ws.Range("J138").Value = 7700.434
$ws.Range("L138").Value = 23101.302
$ws.Range("N138").Value = -33381.302
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4158.5166
$ws.Range("I32").Value = 3125.6023
$ws.Range("J32").Value = 14875.0
$ws.Range("K32").Value = 3125.6023
$ws.Range("L32").Value = 14875.0
$ws.Range("M32").Value = -2838.6023
$ws.Range("N32").Value = -15449.0
$ws.Range("H45").Value = 2172.2083
$ws.Range("I45").Value = 1806.65
$ws.Range("K45").Value = 1806.65
$ws.Range("M45").Value = -1429.65
$ws.Range("H61").Value = 3118.7036
$ws.Range("I61").Value = 3231.875
$ws.Range("J61").Value = 2954.0908
$ws.Range("K61").Value = 3231.875
$ws.Range("L61").Value = 2954.0908
$ws.Range("M61").Value = -3019.875
$ws.Range("N61").Value = -3378.0908
$ws.Range("H97").Value = 625.74194
$ws.Range("I97").Value = 615.0345
$ws.Range("K97").Value = 615.0345
$ws.Range("M97").Value = -119.0345
$ws.Range("H107").Value = 60228.0
$ws.Range("J107").Value = 60228.0
$ws.Range("L107").Value = 60228.0
$ws.Range("N107").Value = -67908.0
$ws.Range("H110").Value = 168195.14
$ws.Range("I110").Value = 209844.05
$ws.Range("J110").Value = 1599.5
$ws.Range("K110").Value = 209844.05
$ws.Range("L110").Value = 1599.5
$ws.Range("M110").Value = -207799.05
$ws.Range("N110").Value = -5689.5
$ws.Range("H111").Value = 67976.664
$ws.Range("J111").Value = 67976.664
$ws.Range("L111").Value = 67976.664
$ws.Range("N111").Value = -76156.664
$ws.Range("H112").Value = 39694.75
$ws.Range("J112").Value = 41259.668
$ws.Range("L112").Value = 41259.668
$ws.Range("N112").Value = -44213.668
$ws.Range("H114").Value = 0.0
$ws.Range("J114").Value = 0.0
$ws.Range("L114").Value = 0.0
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 5791.25
$ws.Range("H136").Value = 3118.7036
$ws.Range("I136").Value = 3231.875
$ws.Range("J136").Value = 2954.0908
$ws.Range("K136").Value = 9695.625
$ws.Range("L136").Value = 8862.2724
$ws.Range("M136").Value = -7145.625
$ws.Range("N136").Value = -13962.2724
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1190.907
$ws.Range("I94").Value = 1216.3334
$ws.Range("K94").Value = 1216.3334
$ws.Range("M94").Value = -765.3334
$ws.Range("H105").Value = 2451.4546
$ws.Range("I105").Value = 2421.6316
$ws.Range("K105").Value = 2421.6316
$ws.Range("M105").Value = -674.6316000000002
$ws.Range("H134").Value = 25402.666
$ws.Range("I134").Value = 3281.0977
$ws.Range("K134").Value = 9843.293099999999
$ws.Range("M134").Value = -7308.293099999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 57742.0
$ws.Range("I31").Value = 1411.2142
$ws.Range("K31").Value = 1411.2142
$ws.Range("M31").Value = -1116.2142
$ws.Range("H34").Value = 57742.0
$ws.Range("I34").Value = 1411.2142
$ws.Range("K34").Value = 1411.2142
$ws.Range("M34").Value = -1209.2142
$ws.Range("H122").Value = 2779.4119
$ws.Range("I122").Value = 2616.8667
$ws.Range("K122").Value = 7850.6001
$ws.Range("M122").Value = -5400.6001
$ws.Range("H134").Value = 359920.06
$ws.Range("I134").Value = 3073.6086
$ws.Range("J134").Value = 2001413.8
$ws.Range("K134").Value = 9220.8258
$ws.Range("L134").Value = 6004241.4
$ws.Range("M134").Value = -6685.825800000001
$ws.Range("N134").Value = -6009311.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0.0
$ws.Range("I80").Value = 0.0
$ws.Range("J80").Value = 0.0
$ws.Range("K80").Value = 0.0
$ws.Range("L80").Value = 0.0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0.0
$ws.Range("I83").Value = 0.0
$ws.Range("J83").Value = 0.0
$ws.Range("K83").Value = 0.0
$ws.Range("L83").Value = 0.0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H131").Value = 2885.8572
$ws.Range("J131").Value = 3142.0
$ws.Range("L131").Value = 9426.0
$ws.Range("N131").Value = -19506.0
$ws.Range("H132").Value = 1005394.0
$ws.Range("I132").Value = 255110.75
$ws.Range("J132").Value = 1434127.2
$ws.Range("K132").Value = 2295996.75
$ws.Range("L132").Value = 12907144.8
$ws.Range("M132").Value = -2293466.75
$ws.Range("N132").Value = -12912204.8
$ws.Range("H136").Value = 2899.0
$ws.Range("I136").Value = 2899.0
$ws.Range("K136").Value = 8697.0
$ws.Range("M136").Value = -3597.0
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 956860.56
$ws.Range("I80").Value = 914617.94
$ws.Range("J80").Value = 1003327.5
$ws.Range("K80").Value = 914617.94
$ws.Range("L80").Value = 1003327.5
$ws.Range("M80").Value = -913619.94
$ws.Range("N80").Value = -1005323.5
$ws.Range("H83").Value = 956860.56
$ws.Range("I83").Value = 914617.94
$ws.Range("J83").Value = 1003327.5
$ws.Range("K83").Value = 4573089.699999999
$ws.Range("L83").Value = 5016637.5
$ws.Range("M83").Value = -4568097.699999999
$ws.Range("N83").Value = -5026621.5
$ws.Range("H93").Value = 49999.0
$ws.Range("J93").Value = 49999.0
$ws.Range("L93").Value = 49999.0
$ws.Range("N93").Value = -53743.0
$ws.Range("H113").Value = 564492.3
$ws.Range("J113").Value = 13449.1
$ws.Range("L113").Value = 13449.1
$ws.Range("N113").Value = -17789.1
$ws.Range("H122").Value = 2846.6875
$ws.Range("I122").Value = 2563.9583
$ws.Range("J122").Value = 3694.875
$ws.Range("K122").Value = 7691.874899999999
$ws.Range("L122").Value = 11084.625
$ws.Range("M122").Value = -5241.874899999999
$ws.Range("N122").Value = -15984.625
$ws.Range("H132").Value = 25138.467
$ws.Range("I132").Value = 3193.3
$ws.Range("K132").Value = 9579.900000000001
$ws.Range("M132").Value = -7049.900000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3653.2122
$ws.Range("I93").Value = 3591.0952
$ws.Range("J93").Value = 3761.9167
$ws.Range("K93").Value = 3591.0952
$ws.Range("L93").Value = 3761.9167
$ws.Range("M93").Value = -2343.0952
$ws.Range("N93").Value = -6257.9167
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 46090.78
$ws.Range("I132").Value = 1668.125
$ws.Range("J132").Value = 147628.28
$ws.Range("K132").Value = 5004.375
$ws.Range("L132").Value = 442884.84
$ws.Range("M132").Value = -2474.375
$ws.Range("N132").Value = -447944.84
